$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("گزارش تردد")
$ws.Rows.Item(3).Delete()
